# P230 preferred names and changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 (CHEONG, DAYEON) preferred name changed from "Dayeon" to "Yeon"
$ws.Range("A9").Value = "Yeon"

# Two new preferred-name entries appended below the existing roster
$ws.Range("A21").Value = "Nora"
$ws.Range("A22").Value = "Allen"

# Match the selection left behind in the saved file
$ws.Range("A22").Select()
